# Auto-generated Excel COM-interop script to update crypto price/volume data
# per the target diff (Wed Feb 28 09:53:05 UTC 2024 GitHub Actions commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string cell updates (Coin names, Links, non-numeric-looking
#     prices, and all Volume% cells) -- assigning these via .Value keeps them
#     stored as text because Excel cannot parse them as a single numeric value.
$ws.Range("D2").Value = "59.282.17"
$ws.Range("E2").Value = "  +4.77%  "
$ws.Range("D3").Value = "3.352.51"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +4.26%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.61%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "3.869.84"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "3.349.53"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "59.042.80"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E33").Value = "  +10.43%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E34").Value = "  +8.52%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  +9.55%  "
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "2.203.45"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  -11.38%  "
$ws.Range("E51").Value = "  +6.57%  "

# --- Price cells whose new text *looks* like a plain number (e.g. "0.998").
#     Excel would normally coerce a bare numeric-looking string into a real
#     number, so temporarily force the cell to Text format, assign the value,
#     then restore the default "Normal" style so no extra formatting lingers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.639"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0981"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "303.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0540"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "137.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.280"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.42"
$ws.Range("D51").Style = "Normal"
